# Bug with the alpha fixed: add the missing data row (ID 14 with its
# date-range string) to Sheet1, then move the selection the way the
# author left it before starting the next tool (plotting several plots
# for one mouse across a range of dates).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row: A3 = 14 (numeric id), B3 = the date-range label that is
# already used elsewhere in the sheet (shared string reused).
$ws.Range("A3").Value = 14
$ws.Range("B3").Value = "(2023-05-4, 2023-5-15)"

# Leave the selection where the author left it.
$ws.Range("B8").Select()
